# The edit swaps a subset of column values between row 11 and row 12
# (the two data rows get new values for columns A, B, E, F, G, H, Q, R, S, Z, AB
#  while the remaining columns stay as-is).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: set to what used to be row 12's values ---
$ws.Range("A11").Value = 111870057
$ws.Range("B11").Value = 90844
$ws.Range("E11").Value = 5449
$ws.Range("F11").Value = "Svart taggsvamp"
$ws.Range("G11").Value = "Phellodon niger"
$ws.Range("H11").Value = "(Fr.:Fr.) P.Karst."
$ws.Range("Q11").Value = 494314
$ws.Range("R11").Value = 6928937
$ws.Range("S11").Value = 20
$ws.Range("Z11").Value = "14:23"
$ws.Range("AB11").Value = "14:23"

# --- Row 12: set to what used to be row 11's values ---
$ws.Range("A12").Value = 111869523
$ws.Range("B12").Value = 56575
$ws.Range("E12").Value = 103021
$ws.Range("F12").Value = "Talltita"
$ws.Range("G12").Value = "Poecile montanus"
$ws.Range("H12").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("Q12").Value = 494333
$ws.Range("R12").Value = 6928943
$ws.Range("S12").Value = 30
$ws.Range("Z12").Value = "15:06"
$ws.Range("AB12").Value = "15:06"
